$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 8000
$ws.Range("J13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("N13").Value = -8338
$ws.Range("H32").Value = 3012.75
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 2683.6667
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 2683.6667
$ws.Range("M32").Value = -3674
$ws.Range("N32").Value = -3335.6667
$ws.Range("H106").Value = 3650
$ws.Range("I106").Value = 3985.7144
$ws.Range("K106").Value = 3985.7144
$ws.Range("M106").Value = -3354.7144
$ws.Range("H138").Value = 2405.0947
$ws.Range("J138").Value = 2303.5645
$ws.Range("L138").Value = 6910.693499999999
$ws.Range("N138").Value = -17190.6935

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4103.34
$ws.Range("I32").Value = 3229.912
$ws.Range("J32").Value = 12934.667
$ws.Range("K32").Value = 3229.912
$ws.Range("L32").Value = 12934.667
$ws.Range("M32").Value = -2942.912
$ws.Range("N32").Value = -13508.667
$ws.Range("H45").Value = 1500.625
$ws.Range("I45").Value = 1277.125
$ws.Range("J45").Value = 1724.125
$ws.Range("K45").Value = 1277.125
$ws.Range("L45").Value = 1724.125
$ws.Range("M45").Value = -900.125
$ws.Range("N45").Value = -2478.125
$ws.Range("H74").Value = 648.881
$ws.Range("I74").Value = 529.075
$ws.Range("K74").Value = 529.075
$ws.Range("M74").Value = 344.925
$ws.Range("H77").Value = 648.881
$ws.Range("I77").Value = 529.075
$ws.Range("K77").Value = 2645.375
$ws.Range("M77").Value = 1722.625
$ws.Range("H109").Value = 68000
$ws.Range("J109").Value = 68000
$ws.Range("L109").Value = 68000
$ws.Range("N109").Value = -70774
$ws.Range("H110").Value = 632.25
$ws.Range("I110").Value = 662.4545000000001
$ws.Range("K110").Value = 662.4545000000001
$ws.Range("M110").Value = 1382.5455
$ws.Range("H123").Value = 82000
$ws.Range("J123").Value = 82000
$ws.Range("L123").Value = 82000
$ws.Range("N123").Value = -91800
$ws.Range("H132").Value = 1753.037
$ws.Range("I132").Value = 1299.8064
$ws.Range("K132").Value = 3899.4192
$ws.Range("M132").Value = -1369.4192

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3089.4285
$ws.Range("I20").Value = 2671.7334
$ws.Range("J20").Value = 4133.6665
$ws.Range("K20").Value = 2671.7334
$ws.Range("L20").Value = 4133.6665
$ws.Range("M20").Value = -2424.7334
$ws.Range("N20").Value = -4627.6665
$ws.Range("H107").Value = 944.36365
$ws.Range("I107").Value = 1036.875
$ws.Range("J107").Value = 697.6667
$ws.Range("K107").Value = 1036.875
$ws.Range("L107").Value = 697.6667
$ws.Range("M107").Value = 883.125
$ws.Range("N107").Value = -4537.6667
$ws.Range("H134").Value = 3607.4639
$ws.Range("I134").Value = 3671.3774
$ws.Range("J134").Value = 3395.75
$ws.Range("K134").Value = 11014.1322
$ws.Range("L134").Value = 10187.25
$ws.Range("M134").Value = -8479.1322
$ws.Range("N134").Value = -15257.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 650.3913
$ws.Range("I16").Value = 597.8333
$ws.Range("K16").Value = 597.8333
$ws.Range("M16").Value = -310.8333
$ws.Range("H31").Value = 1981.2
$ws.Range("J31").Value = 3121.8
$ws.Range("L31").Value = 3121.8
$ws.Range("N31").Value = -3711.8
$ws.Range("H34").Value = 1981.2
$ws.Range("J34").Value = 3121.8
$ws.Range("L34").Value = 3121.8
$ws.Range("N34").Value = -3525.8
$ws.Range("H62").Value = 2562.4
$ws.Range("J62").Value = 2637.3333
$ws.Range("L62").Value = 2637.3333
$ws.Range("N62").Value = -3885.3333
$ws.Range("H65").Value = 2562.4
$ws.Range("J65").Value = 2637.3333
$ws.Range("L65").Value = 13186.6665
$ws.Range("N65").Value = -19426.6665
$ws.Range("H105").Value = 1677.4
$ws.Range("I105").Value = 1677.4
$ws.Range("K105").Value = 1677.4
$ws.Range("M105").Value = 69.59999999999991
$ws.Range("H107").Value = 970.24
$ws.Range("I107").Value = 831.2857
$ws.Range("K107").Value = 831.2857
$ws.Range("M107").Value = 1088.7143
$ws.Range("H113").Value = 650.3913
$ws.Range("I113").Value = 597.8333
$ws.Range("K113").Value = 597.8333
$ws.Range("M113").Value = 1572.1667
$ws.Range("H141").Value = 56640.46
$ws.Range("J141").Value = 55360.5
$ws.Range("L141").Value = 55360.5
$ws.Range("N141").Value = -65720.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 816774.6
$ws.Range("I4").Value = 954695.4399999999
$ws.Range("K4").Value = 2864086.32
$ws.Range("M4").Value = -2863974.32
$ws.Range("H5").Value = 541.4400000000001
$ws.Range("I5").Value = 493.13635
$ws.Range("K5").Value = 1479.40905
$ws.Range("M5").Value = -1367.40905
$ws.Range("H15").Value = 4000
$ws.Range("I15").Value = 4000
$ws.Range("K15").Value = 12000
$ws.Range("M15").Value = -11860
$ws.Range("H121").Value = 674.1539
$ws.Range("I121").Value = 571.8570999999999
$ws.Range("J121").Value = 793.5
$ws.Range("K121").Value = 1715.5713
$ws.Range("L121").Value = 2380.5
$ws.Range("M121").Value = -405.5712999999998
$ws.Range("N121").Value = -5000.5
$ws.Range("H122").Value = 924.1923
$ws.Range("I122").Value = 599.6667
$ws.Range("K122").Value = 5397.0003
$ws.Range("M122").Value = -2947.0003
$ws.Range("H131").Value = 26535.625
$ws.Range("J131").Value = 29242.828
$ws.Range("L131").Value = 87728.484
$ws.Range("N131").Value = -97808.484
$ws.Range("H135").Value = 541.4400000000001
$ws.Range("I135").Value = 493.13635
$ws.Range("K135").Value = 4438.22715
$ws.Range("M135").Value = -1903.22715

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 64312.223
$ws.Range("J110").Value = 64312.223
$ws.Range("L110").Value = 64312.223
$ws.Range("N110").Value = -72492.223
$ws.Range("H132").Value = 1204293.4
$ws.Range("I132").Value = 1833284
$ws.Range("K132").Value = 5499852
$ws.Range("M132").Value = -5497322

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2374.7827
$ws.Range("I61").Value = 2255.4546
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2255.4546
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2053.4546
$ws.Range("N61").Value = -5404
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H113").Value = 2374.7827
$ws.Range("I113").Value = 2255.4546
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2255.4546
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -85.45460000000003
$ws.Range("N113").Value = -9340
$ws.Range("H132").Value = 5900.737
$ws.Range("I132").Value = 5210.7144
$ws.Range("K132").Value = 15632.1432
$ws.Range("M132").Value = -13102.1432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 10867.889
$ws.Range("J96").Value = 13030.571
$ws.Range("L96").Value = 13030.571
$ws.Range("N96").Value = -15776.571
$ws.Range("H107").Value = 594.6
$ws.Range("I107").Value = 418.7
$ws.Range("K107").Value = 1256.1
$ws.Range("M107").Value = 663.9000000000001
$ws.Range("H113").Value = 1124.1428
$ws.Range("J113").Value = 1079.625
$ws.Range("L113").Value = 3238.875
$ws.Range("N113").Value = -7578.875
$ws.Range("H123").Value = 68000
$ws.Range("J123").Value = 68000
$ws.Range("L123").Value = 68000
$ws.Range("N123").Value = -77800
